# Update references slightly for better accuracy.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 29

# --- Row 4 ---
$ws.Range("B4").Value = 15
$ws.Range("C4").Value = 52

# --- Row 5 --- (also drop the "Ideal" note in H5)
$ws.Range("F5").Value = 9
$ws.Range("G5").Value = 26
$ws.Range("H5").ClearContents()

# --- Row 6 ---
$ws.Range("B6").Value = 12
$ws.Range("C6").Value = 49
$ws.Range("D6").Value = 44
$ws.Range("E6").Value = 68
$ws.Range("F6").Value = 10
$ws.Range("G6").Value = 28

# --- Row 7 ---
$ws.Range("B7").Value = 9
$ws.Range("C7").Value = 46
$ws.Range("D7").Value = 37
$ws.Range("E7").Value = 61
$ws.Range("F7").Value = 9

# --- Row 8 ---
$ws.Range("B8").Value = 11
$ws.Range("C8").Value = 48
$ws.Range("D8").Value = 29
$ws.Range("E8").Value = 53
$ws.Range("G8").Value = 26

# --- Row 9 ---
$ws.Range("B9").Value = 10
$ws.Range("C9").Value = 47
$ws.Range("D9").Value = 37
$ws.Range("E9").Value = 61
$ws.Range("F9").Value = 13
$ws.Range("G9").Value = 30

# --- Rows 10 & 11 swap places (sorted by Accession Num), row 10's
#     values also change slightly; the "Smaller by around 10-20%" note
#     that used to live on H11 is dropped. ---
$ws.Range("A10").Value = 8649024
$ws.Range("B10").Value = 22
$ws.Range("C10").Value = 59
$ws.Range("D10").Value = 21
$ws.Range("E10").Value = 45
$ws.Range("F10").Value = 11
$ws.Range("G10").Value = 29

$ws.Range("A11").Value = 8931305
$ws.Range("B11").Value = 19
$ws.Range("C11").Value = 56
$ws.Range("D11").Value = 46
$ws.Range("E11").Value = 70
$ws.Range("F11").Value = 13
$ws.Range("G11").Value = 37
$ws.Range("H11").ClearContents()

# --- Row 12 ---
$ws.Range("D12").Value = 38
$ws.Range("E12").Value = 62

# C10 keeps an integer ("0") number format, same as the rest of the
# B:G data columns.
$ws.Range("C10").NumberFormat = "0"

# Move the active selection to H5 (matches the saved cursor position).
$ws.Range("H5").Select()
